# Update the "Förändrad" date column (column C) for all data rows.
# The workbook stores dates as serial numbers; every value of 45180
# (2023-09-11) in C2:C79 is bumped by one day to 45181 (2023-09-12).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row  # xlUp = -4162
if ($lastRow -lt 2) { $lastRow = 2 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45180) {
        $cell.Value = 45181
    }
}
